# NBA DATA Changes for recomadation end date
#
# The "recommendation_end_date" column (AE) on the data_NBA sheet holds the
# value "2019-03-31" for almost every data row. A handful of rows (6, 29, 30
# and 31) had stray / inconsistent end-date values. This change normalises
# those four rows to the same end date as the rest of the sheet, and at the
# same time rolls the recommendation end date forward from 2019-03-31 to
# 2021-03-31 for every row that has one (rows whose end date is "NULL" are
# left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldEndDate = "2019-03-31"
$newEndDate = "2021-03-31"

$lastRow = 56
for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Range("AE" + $row)
    if ($cell.Text -ne "NULL") {
        $cell.Value = $newEndDate
    }
}
